# Updates the cryptocurrency price/label snapshot on Sheet1 to match the
# latest scrape (commit: "Updated symbol list on Sat Dec 24 07:17:42 UTC 2022
# with GitHub Actions").
#
# The "Price" column (D) holds numeric-looking values that are stored as text
# in the workbook, so each is written with a leading apostrophe to force
# Excel to keep them as text instead of silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$Value
    )
    $ws.Range($CellRef).Value = "'" + $Value
}

# Price column (D) updates
Set-TextValue "D2"  "245.60"
Set-TextValue "D3"  "22.04"
Set-TextValue "D4"  "5.356"
Set-TextValue "D5"  "0.05963"
Set-TextValue "D7"  "6.392"
Set-TextValue "D8"  "0.8101"
Set-TextValue "D9"  "0.9632"
Set-TextValue "D10" "0.1430"
Set-TextValue "D11" "0.07383"
Set-TextValue "D14" "0.09407"
Set-TextValue "D16" "0.001597"
Set-TextValue "D17" "0.04788"
Set-TextValue "D19" "0.006211"
Set-TextValue "D20" "0.005105"
Set-TextValue "D21" "0.0009841"
Set-TextValue "D22" "0.00006903"
Set-TextValue "D23" "3.747"
Set-TextValue "D40" "0.04033"
Set-TextValue "D41" "0.006453"
Set-TextValue "D43" "0.002721"
Set-TextValue "D44" "0.005804"
Set-TextValue "D45" "0.00005258"
Set-TextValue "D48" "0.03295"

# Volume(1h) column (E) label updates
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
